$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported for "Región de Arica y Parinacota".
# Insert it as a new row 19 (pushing the existing rows 19-53 down to 20-54,
# which is exactly what the diff shows: every old row N's data reappears at N+1).
$ws.Rows("19:19").Insert()

$ws.Range("A19").Value = 11
$ws.Range("B19").Value = "Vega Monumental Concepción"
$ws.Range("C19").Value = "Bíobío"
$ws.Range("D19").Value = 44533
$ws.Range("E19").Value = 8
$ws.Range("F19").Value = 100112001
$ws.Range("G19").Value = "Berenjena"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 180
$ws.Range("K19").Value = 7500
$ws.Range("L19").Value = 8000
$ws.Range("M19").Value = 7778
$ws.Range("N19").Value = "$/caja 60 unidades"
$ws.Range("O19").Value = "Región de Arica y Parinacota"
$ws.Range("P19").Value = 130
$ws.Range("Q19").Value = 60
$ws.Range("R19").Value = "Hortaliza"
